{"js": "// Update the division-problem table: 25 cells (5 data rows x 5 cols,\n// stored as OOXML rows 0,4,8,12,16) get new \"a\u00f7b=q, r\" values.\n// Net cell/row/column counts are unchanged (an insertion of one cell\n// and a deletion of one cell in the 3rd data row cancel out), so every\n// cell is addressed directly via table.getCell(row, col) and its\n// .value is replaced in place -- this preserves the existing run\n// formatting (TimeNewRoman, sz 30, left-justified paragraph) already\n// present in each cell.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Map of absolute OOXML row index -> new texts for columns 0..4.\nconst updates = {\n  0: [\"28\u00f72=14, 0\", \"77\u00f76=12, 5\", \"14\u00f75=2, 4\", \"83\u00f74=20, 3\", \"72\u00f76=12, 0\"],\n  4: [\"35\u00f76=5, 5\", \"12\u00f73=4, 0\", \"26\u00f74=6, 2\", \"11\u00f79=1, 2\", \"81\u00f75=16, 1\"],\n  8: [\"43\u00f73=14, 1\", \"74\u00f76=12, 2\", \"89\u00f78=11, 1\", \"57\u00f77=8, 1\", \"62\u00f79=6, 8\"],\n  12: [\"28\u00f74=7, 0\", \"11\u00f73=3, 2\", \"85\u00f78=10, 5\", \"68\u00f79=7, 5\", \"24\u00f75=4, 4\"],\n  16: [\"65\u00f77=9, 2\", \"71\u00f76=11, 5\", \"73\u00f75=14, 3\", \"66\u00f73=22, 0\", \"44\u00f73=14, 2\"],\n};\n\nfor (const rowIndexStr of Object.keys(updates)) {\n  const rowIndex = parseInt(rowIndexStr, 10);\n  const values = updates[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(rowIndex, col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the division-problem table: 25 cells (5 data rows x 5 cols,\n# stored as table rows 1,5,9,13,17 in 1-based COM indexing) get new\n# \"a\u00f7b=q, r\" values. The net cell/row/column counts are unchanged (an\n# insertion of one cell and a deletion of one cell in the 3rd data row\n# cancel out), so every cell is addressed directly via $table.Cell(r,c)\n# and its Range.Text is replaced in place -- this preserves the\n# existing run formatting (TimeNewRoman, sz 30, left-justified\n# paragraph) already present in each cell.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$updates = @{\n  1  = @(\"28\u00f72=14, 0\", \"77\u00f76=12, 5\", \"14\u00f75=2, 4\", \"83\u00f74=20, 3\", \"72\u00f76=12, 0\")\n  5  = @(\"35\u00f76=5, 5\", \"12\u00f73=4, 0\", \"26\u00f74=6, 2\", \"11\u00f79=1, 2\", \"81\u00f75=16, 1\")\n  9  = @(\"43\u00f73=14, 1\", \"74\u00f76=12, 2\", \"89\u00f78=11, 1\", \"57\u00f77=8, 1\", \"62\u00f79=6, 8\")\n  13 = @(\"28\u00f74=7, 0\", \"11\u00f73=3, 2\", \"85\u00f78=10, 5\", \"68\u00f79=7, 5\", \"24\u00f75=4, 4\")\n  17 = @(\"65\u00f77=9, 2\", \"71\u00f76=11, 5\", \"73\u00f75=14, 3\", \"66\u00f73=22, 0\", \"44\u00f73=14, 2\")\n}\n\nforeach ($rowIndex in $updates.Keys) {\n  $values = $updates[$rowIndex]\n  for ($col = 1; $col -le $values.Length; $col++) {\n    $cell = $table.Cell($rowIndex, $col)\n    $cell.Range.Text = $values[$col - 1]\n  }\n}\n"}
